$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("B2").Value = 9.573320881146669
$ws.Range("C2").Value = 6.069437374781471
$ws.Range("E2").Value = 16.45082476823774
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.630909019046887
$ws.Range("K2").Value = 8.733770021307958
$ws.Range("O2").Value = 21.25343720909476

$ws.Range("B3").Value = 9.226402674923111
$ws.Range("C3").Value = 5.906420148546537
$ws.Range("E3").Value = 15.52155248879117
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.632843367313531
$ws.Range("K3").Value = 8.481627826679846
$ws.Range("O3").Value = 21.37703864509586

$ws.Range("B4").Value = 9.00793307964571
$ws.Range("C4").Value = 5.803314188576942
$ws.Range("E4").Value = 14.92639962209146
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.634091431664017
$ws.Range("K4").Value = 8.323990692138057
$ws.Range("O4").Value = 21.45872936310025

$ws.Range("B5").Value = 8.917669344557927
$ws.Range("C5").Value = 5.760581508676516
$ws.Range("E5").Value = 14.67796088590094
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.634615259677332
$ws.Range("K5").Value = 8.259137829946045
$ws.Range("O5").Value = 21.49347214487979

$ws.Range("B6").Value = 8.90261077481008
$ws.Range("C6").Value = 5.753443796061382
$ws.Range("E6").Value = 14.63635988220539
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.634703162401581
$ws.Range("K6").Value = 8.248334845727094
$ws.Range("O6").Value = 21.4993287785529

$ws.Range("B7").Value = 9.006720562707326
$ws.Range("C7").Value = 5.802740723423362
$ws.Range("E7").Value = 14.92307261963871
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.634098434450983
$ws.Range("K7").Value = 8.32311842168609
$ws.Range("O7").Value = 21.45919203931047

$ws.Range("B8").Value = 9.454915107218502
$ws.Range("C8").Value = 6.013878206343492
$ws.Range("E8").Value = 16.13565920298397
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.631563484578738
$ws.Range("K8").Value = 8.647471792028469
$ws.Range("O8").Value = 21.29484836290758

$ws.Range("B9").Value = 10.28513461526957
$ws.Range("C9").Value = 6.402353382730895
$ws.Range("E9").Value = 18.37948890688392
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.627069057333558
$ws.Range("K9").Value = 9.257307427177839
$ws.Range("O9").Value = 21.01881125754996

$ws.Range("B10").Value = 10.85914066461815
$ws.Range("C10").Value = 6.670178003545413
$ws.Range("E10").Value = 20.01175326496732
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.624054204228157
$ws.Range("K10").Value = 9.684702962208888
$ws.Range("O10").Value = 20.84451371343301

$ws.Range("B11").Value = 11.11137926513759
$ws.Range("C11").Value = 6.787836018630869
$ws.Range("E11").Value = 20.71227790280468
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.622744318217672
$ws.Range("K11").Value = 9.873800923481076
$ws.Range("O11").Value = 20.77148172374639

$ws.Range("B12").Value = 11.20554288335132
$ws.Range("C12").Value = 6.831762894332279
$ws.Range("E12").Value = 20.97154488033283
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.622257100142191
$ws.Range("K12").Value = 9.944579401789088
$ws.Range("O12").Value = 20.74473224699331

$ws.Range("B13").Value = 11.18532438615326
$ws.Range("C13").Value = 6.822330773085331
$ws.Range("E13").Value = 20.91597371328208
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.622361640335531
$ws.Range("K13").Value = 9.929373786240269
$ws.Range("O13").Value = 20.7504528204653

$ws.Range("B14").Value = 11.11915369693935
$ws.Range("C14").Value = 6.7914626337152
$ws.Range("E14").Value = 20.73372814329073
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.622704058259583
$ws.Range("K14").Value = 9.879640859972399
$ws.Range("O14").Value = 20.76926283032824

$ws.Range("B15").Value = 11.07844384188745
$ws.Range("C15").Value = 6.772472510531233
$ws.Range("E15").Value = 20.62131622543686
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.622914944660165
$ws.Range("K15").Value = 9.84906831186675
$ws.Range("O15").Value = 20.78090270502018

$ws.Range("B16").Value = 10.84247080578436
$ws.Range("C16").Value = 6.662402354090066
$ws.Range("E16").Value = 19.96513001069957
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.624141043459503
$ws.Range("K16").Value = 9.672232264109807
$ws.Range("O16").Value = 20.84941297607295

$ws.Range("B17").Value = 10.69537686643565
$ws.Range("C17").Value = 6.593787995493818
$ws.Range("E17").Value = 19.55184730954957
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.624908954269361
$ws.Range("K17").Value = 9.562338085608499
$ws.Range("O17").Value = 20.89304887119923

$ws.Range("B18").Value = 10.6099404521799
$ws.Range("C18").Value = 6.553931134504262
$ws.Range("E18").Value = 19.31018276877987
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.625356435890069
$ws.Range("K18").Value = 9.498631859131049
$ws.Range("O18").Value = 20.91873561526581

$ws.Range("B19").Value = 10.58087278001999
$ws.Range("C19").Value = 6.540369880951618
$ws.Range("E19").Value = 19.22767897404509
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.625508943058352
$ws.Range("K19").Value = 9.476978513689636
$ws.Range("O19").Value = 20.92753359239071

$ws.Range("B20").Value = 10.71112195068495
$ws.Range("C20").Value = 6.601132854833964
$ws.Range("E20").Value = 19.59625107267533
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.624826608964726
$ws.Range("K20").Value = 9.574088525428138
$ws.Range("O20").Value = 20.88834279862811

$ws.Range("B21").Value = 11.13862693353154
$ws.Range("C21").Value = 6.800546591204827
$ws.Range("E21").Value = 20.78742083661344
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.622603243172413
$ws.Range("K21").Value = 9.894271589674917
$ws.Range("O21").Value = 20.76371322995483

$ws.Range("B22").Value = 11.4101071075177
$ws.Range("C22").Value = 6.927205836542033
$ws.Range("E22").Value = 21.53092387598497
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.621201458571746
$ws.Range("K22").Value = 10.09867397227307
$ws.Range("O22").Value = 20.68754541631431

$ws.Range("B23").Value = 11.26596048341411
$ws.Range("C23").Value = 6.859949386043132
$ws.Range("E23").Value = 21.1372934906544
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.621944938208321
$ws.Range("K23").Value = 9.990044036539166
$ws.Range("O23").Value = 20.72771189830684

$ws.Range("B24").Value = 10.70400630246914
$ws.Range("C24").Value = 6.597813517568812
$ws.Range("E24").Value = 19.57618878926483
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.624863818584286
$ws.Range("K24").Value = 9.568777786961375
$ws.Range("O24").Value = 20.89046854579862

$ws.Range("B25").Value = 10.06644666523732
$ws.Range("C25").Value = 6.300216107374585
$ws.Range("E25").Value = 17.74414544876901
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.628234242450706
$ws.Range("K25").Value = 9.095645975827601
$ws.Range("O25").Value = 21.08850390262307
